# B6-PowerPoint.pptx edit:
#   1. Three tables (slides 14, 15, 16) switch from the deck's single
#      custom table style {CC85F6E7-CD32-44B5-BED5-79E9C7031EF1} to the
#      built-in style {D32A397A-1999-4E44-B5A6-C4F98FDC5BF0}.
#   2. The two theme parts swap identity: the theme actually in force for
#      the slide master/slides (formerly "Integral" / Red-Violet colours)
#      becomes the "Office Theme" colour set. (Font scheme and format
#      scheme are already identical between the two themes, so only the
#      colour scheme actually changes visibly.)

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------
$newStyleId = "{D32A397A-1999-4E44-B5A6-C4F98FDC5BF0}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $tableShape = $slide.Shapes.Item(1)
    $tableShape.Table.ApplyStyle($newStyleId)
}

# --- 2. Swap the colour scheme so the master's theme carries the -----
#        "Office Theme" colours (previously the "Integral" colours)
$theme = $p.SlideMaster.Theme.ThemeColorScheme

# index : 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3
#         8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le 12; $i++) {
    $theme.Colors($i).RGB = $officeColors[$i - 1]
}
